$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

$ws.Range("G3").Value = 2000
$ws.Range("G5").Value = 2500
$ws.Range("G12").Value = 1000
$ws.Range("G14").Value = 3000
$ws.Range("G17").Value = 3500
$ws.Range("G18").Value = 7500
$ws.Range("G19").Value = 5000
$ws.Range("G21").Value = 3000
$ws.Range("G29").Value = 32000
